$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, move the existing "White River Lower 08" row (currently row 2) down to row 4,
# writing directly (no Insert/shift) so no formatting is copied from neighboring rows.
$ws.Range("A4").Value2 = "White River Lower 08"
$ws.Range("B4").Value2 = "Wenatchee"
$ws.Range("C4").Value2 = "Lower White River"
$ws.Range("D4").Value2 = "yes"
$ws.Range("E4").Value2 = "yes"
$ws.Range("F4").Value2 = "yes"
$ws.Range("G4").Value2 = 5
$ws.Range("H4").Value2 = 5
$ws.Range("I4").Value2 = 5
$ws.Range("J4").Value2 = 5
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 5
$ws.Range("M4").Value2 = 5
$ws.Range("N4").Value2 = 1
$ws.Range("O4").Value2 = 1
$ws.Range("P4").Value2 = 3
$ws.Range("Q4").Value2 = 5
$ws.Range("R4").Value2 = 4
$ws.Range("S4").Value2 = 5
$ws.Range("T4").Value2 = 34
$ws.Range("U4").Value2 = 0.7555555555555555
$ws.Range("V4").Value2 = 5
$ws.Range("W4").Value2 = 3

# Now overwrite row 2 with "Twisp River Upper 02" data
$ws.Range("A2").Value2 = "Twisp River Upper 02"
$ws.Range("B2").Value2 = "Methow"
$ws.Range("C2").Value2 = "Upper Twisp River"
$ws.Range("D2").Value2 = "yes"
$ws.Range("E2").Value2 = "yes"
$ws.Range("F2").Value2 = "yes"
$ws.Range("G2").Value2 = 5
$ws.Range("H2").Value2 = 5
$ws.Range("I2").Value2 = 5
$ws.Range("J2").Value2 = 5
$ws.Range("K2").Value2 = 1
$ws.Range("L2").Value2 = 5
$ws.Range("M2").Value2 = 5
$ws.Range("N2").Value2 = 5
$ws.Range("O2").Value2 = 3
$ws.Range("P2").Value2 = 1
$ws.Range("Q2").Value2 = 5
$ws.Range("R2").Value2 = 3
$ws.Range("S2").Value2 = 5
$ws.Range("T2").Value2 = 37
$ws.Range("U2").Value2 = 0.8222222222222222
$ws.Range("V2").Value2 = 3
$ws.Range("W2").Value2 = 3

# Add new row 3 "Twisp River Upper 03" data (same values as row 2, different name)
$ws.Range("A3").Value2 = "Twisp River Upper 03"
$ws.Range("B3").Value2 = "Methow"
$ws.Range("C3").Value2 = "Upper Twisp River"
$ws.Range("D3").Value2 = "yes"
$ws.Range("E3").Value2 = "yes"
$ws.Range("F3").Value2 = "yes"
$ws.Range("G3").Value2 = 5
$ws.Range("H3").Value2 = 5
$ws.Range("I3").Value2 = 5
$ws.Range("J3").Value2 = 5
$ws.Range("K3").Value2 = 1
$ws.Range("L3").Value2 = 5
$ws.Range("M3").Value2 = 5
$ws.Range("N3").Value2 = 5
$ws.Range("O3").Value2 = 3
$ws.Range("P3").Value2 = 1
$ws.Range("Q3").Value2 = 5
$ws.Range("R3").Value2 = 3
$ws.Range("S3").Value2 = 5
$ws.Range("T3").Value2 = 37
$ws.Range("U3").Value2 = 0.8222222222222222
$ws.Range("V3").Value2 = 3
$ws.Range("W3").Value2 = 3
